$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateFormat = $ws.Range("A51").NumberFormat

$rows = @(
    @{ Row = 52; Date = 43810.45908564814; User = "jetnew"; Text = "test"; Answer = "No content found." },
    @{ Row = 53; Date = 43810.45924768518; User = "jetnew"; Text = "What does it mean to be curious?"; Answer = "No content found." },
    @{ Row = 54; Date = 43810.45946759259; User = "jetnew"; Text = "how long do students live in cinnamon college?"; Answer = "No content found." }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Date
    $ws.Cells.Item($r.Row, 1).NumberFormat = $dateFormat
    $ws.Cells.Item($r.Row, 2).Value = $r.User
    $ws.Cells.Item($r.Row, 3).Value = $r.Text
    $ws.Cells.Item($r.Row, 4).Value = $r.Answer
}
